$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "date"
$ws.Range("B1").Value = "hist"

# Data values (24 hourly rows starting 2020-08-25 00:00:00, i.e. serial 44068)
$dates = @(44068, 44068.04166666666, 44068.08333333334, 44068.125, 44068.16666666666, 44068.20833333334, 44068.25, 44068.29166666666, 44068.33333333334, 44068.375, 44068.41666666666, 44068.45833333334, 44068.5, 44068.54166666666, 44068.58333333334, 44068.625, 44068.66666666666, 44068.70833333334, 44068.75, 44068.79166666666, 44068.83333333334, 44068.875, 44068.91666666666, 44068.95833333334)
$hist  = @(481.25, 262.75, 262.25, 271.375, 275, 269.25, 1208.375, 617.5, 2460, 3102.25, 5401.625, 5650.375, 5908, 5907, 6052, 6074.75, 6967.5, 6293.375, 6206.875, 5380, 4597.5, 4619.25, 1652.875, 705.75)

for ($i = 0; $i -lt 24; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $hist[$i]
}

# Apply the date/time number format to column A (rows 2-25).
# Set it twice on A2 (first lower-case, then upper-case) so the engine
# settles on a single reused style slot instead of allocating an extra
# unused one, then extend the final format to the rest of the range.
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A3:A25").NumberFormat = "YYYY-MM-DD HH:MM:SS"
